# Apply "More hyperparamter updated results" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 53-57: add newly-computed IOU2 / Thresh2 (and helper) columns F..I
# ---------------------------------------------------------------------
$ws.Range("F53").Value = 0.244
$ws.Range("G53").Value = 0.58
$ws.Range("H53").Value = 0.279
$ws.Range("I53").Value = 0.54

$ws.Range("F54").Value = 0.259
$ws.Range("G54").Value = 0.58
$ws.Range("H54").Value = 0.298
$ws.Range("I54").Value = 0.58

$ws.Range("F55").Value = 0.111
$ws.Range("G55").Value = 0.33
$ws.Range("H55").Value = 0.176
$ws.Range("I55").Value = 0.29

$ws.Range("F56").Value = 0.148
$ws.Range("G56").Value = 0.38
$ws.Range("H56").Value = 0.202
$ws.Range("I56").Value = 0.38

$ws.Range("F57").Value = 0.253
$ws.Range("G57").Value = 0.42
$ws.Range("H57").Value = 0.307
$ws.Range("I57").Value = 0.38

# ---------------------------------------------------------------------
# Rows 66-71 (HYPERPARAMETER SEARCH 3 block): updated / new F..I values
# ---------------------------------------------------------------------
$ws.Range("F66").Value = 0
$ws.Range("G66").Value = 0
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0

$ws.Range("F67").Value = 0.038
$ws.Range("G67").Value = 1
$ws.Range("H67").Value = 0.078
$ws.Range("I67").Value = 1

$ws.Range("F68").Value = 0.211
$ws.Range("G68").Value = 0.38
$ws.Range("H68").Value = 0.256
$ws.Range("I68").Value = 0.29

$ws.Range("F69").Value = 0.266
$ws.Range("G69").Value = 0.58
$ws.Range("H69").Value = 0.292
$ws.Range("I69").Value = 0.58

$ws.Range("F70").Value = 0.235
$ws.Range("G70").Value = 0.83
$ws.Range("H70").Value = 0.29
$ws.Range("I70").Value = 0.83

$ws.Range("F71").Value = 0.01
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 0.02
$ws.Range("I71").Value = 0

# ---------------------------------------------------------------------
# Row 74 header: label the previously-blank IOU2 / Thresh2 columns
# ---------------------------------------------------------------------
$ws.Range("H74").Value = "IOU2"
$ws.Range("I74").Value = "Thresh2"

# ---------------------------------------------------------------------
# Row 75: updated F..I values
# ---------------------------------------------------------------------
$ws.Range("F75").Value = 0.368
$ws.Range("G75").Value = 0.83
$ws.Range("H75").Value = 0.406
$ws.Range("I75").Value = 0.83

# ---------------------------------------------------------------------
# New row 76: another run ("run117") inserted right after row 75
# ---------------------------------------------------------------------
$ws.Range("A76").Value = 0.000003
$ws.Range("B76").Value = 0.0005
$ws.Range("C76").Value = "run117"
$ws.Range("D76").Value = 1

# ---------------------------------------------------------------------
# Rows 132-133: add new F..I values
# ---------------------------------------------------------------------
$ws.Range("F132").Value = 0.211
$ws.Range("G132").Value = 0.5
$ws.Range("H132").Value = 0.255
$ws.Range("I132").Value = 0.5

$ws.Range("F133").Value = 0.311
$ws.Range("G133").Value = 0.92
$ws.Range("H133").Value = 0.368
$ws.Range("I133").Value = 0.92

# ---------------------------------------------------------------------
# View-state bookkeeping: scroll position + selection moved down to the
# newly-edited area, and the sheet-tab area was shrunk.
# ---------------------------------------------------------------------
try {
    $excel.ActiveWindow.ScrollRow = 125
    $excel.ActiveWindow.ScrollColumn = 1
} catch {}
try {
    $excel.ActiveWindow.TabRatio = 174
} catch {}

$ws.Range("F134").Select()
